$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing content
$ws.Cells.Clear()

# Header row
$ws.Cells.Item(1,1).Value = "项目名"
$ws.Cells.Item(1,2).Value = "项目链接"
$ws.Cells.Item(1,3).Value = "学院"
$ws.Cells.Item(1,4).Value = "项目简介"
$ws.Cells.Item(1,5).Value = "链接"

# Row 2
$ws.Cells.Item(2,1).Value = "advanced aesthetic dentistry pg cert"
$ws.Cells.Item(2,2).Value = "https://www.ucl.ac.uk/prospective-students/graduate/taught-degrees/advanced-aesthetic-dentistry-pg-cert"
$ws.Cells.Item(2,3).Value = "Faculty of Medical Sciences | Eastman Dental Institute"
$ws.Cells.Item(2,4).Value = "This programme explores the most important topics in aesthetic dentistry, making it ideal for practitioners looking to expand into this sought-after field. It covers scientific documentation, clinical applications, treatment planning, modern materials, and step-by-step clinical procedures. This one-year course is available part-time with regular attendance or by distance-learning with three condensed week-long sessions in London."
$ws.Cells.Item(2,5).Value = "English language requirements:https://www.ucl.ac.uk/prospective-students/graduate/learning-and-living-ucl/international-students/english-language-requirements"
$ws.Cells.Item(2,6).Value = "Application fees:https://www.ucl.ac.uk/prospective-students/graduate/application-fees"
$ws.Cells.Item(2,7).Value = "Entry requirements:#entry-requirements"
$ws.Cells.Item(2,8).Value = "Graduate degrees:/prospective-students/graduate/graduate-degrees"
$ws.Cells.Item(2,9).Value = "Taught Degrees:/prospective-students/graduate/taught-degrees"
$ws.Cells.Item(2,10).Value = "Applying for Graduate Taught Study at UCL:/prospective-students/graduate/taught-degrees/applying-graduate-taught-study-ucl"
$ws.Cells.Item(2,11).Value = "Research Degrees:/prospective-students/graduate/research-degrees"
$ws.Cells.Item(2,12).Value = "Applying for Graduate Research Study at UCL:/prospective-students/graduate/research-degrees/applying-graduate-research-study-ucl"
$ws.Cells.Item(2,13).Value = "Entry requirements:/prospective-students/graduate/teacher-training/entry-requirements"
$ws.Cells.Item(2,14).Value = "How to apply:/prospective-students/graduate/teacher-training/how-apply"

# Row 3
$ws.Cells.Item(3,1).Value = "advanced architectural research pg cert"
$ws.Cells.Item(3,2).Value = "https://www.ucl.ac.uk/prospective-students/graduate/taught-degrees/advanced-architectural-research-pg-cert"
$ws.Cells.Item(3,3).Value = "Faculty of the Built Environment | Bartlett School of Architecture"
$ws.Cells.Item(3,4).Value = "In architecture, research and practice go hand in hand, with ground-breaking design inspiring great research and vice versa. This programme allows students to develop their own research and/or design to an advanced level, ready to push new boundaries and change normalised modes of working either in academia or industry. Many students on this programme work with their tutors to clarify a PhD proposal."
$ws.Cells.Item(3,5).Value = "English language requirements:https://www.ucl.ac.uk/prospective-students/graduate/learning-and-living-ucl/international-students/english-language-requirements"
$ws.Cells.Item(3,6).Value = "Application fees:https://www.ucl.ac.uk/prospective-students/graduate/application-fees"
$ws.Cells.Item(3,7).Value = "Entry requirements:#entry-requirements"
$ws.Cells.Item(3,8).Value = "Graduate degrees:/prospective-students/graduate/graduate-degrees"
$ws.Cells.Item(3,9).Value = "Taught Degrees:/prospective-students/graduate/taught-degrees"
$ws.Cells.Item(3,10).Value = "Applying for Graduate Taught Study at UCL:/prospective-students/graduate/taught-degrees/applying-graduate-taught-study-ucl"
$ws.Cells.Item(3,11).Value = "Research Degrees:/prospective-students/graduate/research-degrees"
$ws.Cells.Item(3,12).Value = "Applying for Graduate Research Study at UCL:/prospective-students/graduate/research-degrees/applying-graduate-research-study-ucl"
$ws.Cells.Item(3,13).Value = "Entry requirements:/prospective-students/graduate/teacher-training/entry-requirements"
$ws.Cells.Item(3,14).Value = "How to apply:/prospective-students/graduate/teacher-training/how-apply"

# Row 4
$ws.Cells.Item(4,1).Value = "advanced audiology msc"
$ws.Cells.Item(4,2).Value = "https://www.ucl.ac.uk/prospective-students/graduate/taught-degrees/advanced-audiology-msc"
$ws.Cells.Item(4,3).Value = "Faculty of Brain Sciences | Ear Institute"
$ws.Cells.Item(4,4).Value = "Advanced Audiology is designed for practising audiologists looking to enhance their clinical skills. Graduates are eligible to apply for senior clinical roles in the NHS or private sector. A variety of specialist modules can be selected to suit your professional needs. The programme has a strong research and evidence-based practice foundation."
$ws.Cells.Item(4,5).Value = "English language requirements:https://www.ucl.ac.uk/prospective-students/graduate/learning-and-living-ucl/international-students/english-language-requirements"
$ws.Cells.Item(4,6).Value = "UCL Graduate Admissions team:https://www.ucl.ac.uk/prospective-students/graduate/admissions-enquiries#form"
$ws.Cells.Item(4,7).Value = "Application fees:https://www.ucl.ac.uk/prospective-students/graduate/application-fees"
$ws.Cells.Item(4,8).Value = "Entry requirements:#entry-requirements"
$ws.Cells.Item(4,9).Value = "Graduate degrees:/prospective-students/graduate/graduate-degrees"
$ws.Cells.Item(4,10).Value = "Taught Degrees:/prospective-students/graduate/taught-degrees"
$ws.Cells.Item(4,11).Value = "Applying for Graduate Taught Study at UCL:/prospective-students/graduate/taught-degrees/applying-graduate-taught-study-ucl"
$ws.Cells.Item(4,12).Value = "Research Degrees:/prospective-students/graduate/research-degrees"
$ws.Cells.Item(4,13).Value = "Applying for Graduate Research Study at UCL:/prospective-students/graduate/research-degrees/applying-graduate-research-study-ucl"
$ws.Cells.Item(4,14).Value = "Entry requirements:/prospective-students/graduate/teacher-training/entry-requirements"
$ws.Cells.Item(4,15).Value = "How to apply:/prospective-students/graduate/teacher-training/how-apply"

# Row 5
$ws.Cells.Item(5,1).Value = "advanced audiology: audiovestibular medicine msc"
$ws.Cells.Item(5,2).Value = "https://www.ucl.ac.uk/prospective-students/graduate/taught-degrees/advanced-audiology-audiovestibular-medicine-msc"
$ws.Cells.Item(5,3).Value = "Faculty of Brain Sciences | Ear Institute"
$ws.Cells.Item(5,4).Value = "The programme aims to provide the theoretical, clinical and research skills necessary for physicians who wish to pursue a career in audiovestibular medicine."
$ws.Cells.Item(5,5).Value = "English language requirements:https://www.ucl.ac.uk/prospective-students/graduate/learning-and-living-ucl/international-students/english-language-requirements"
$ws.Cells.Item(5,6).Value = "UCL Graduate Admissions team:https://www.ucl.ac.uk/prospective-students/graduate/admissions-enquiries#form"
$ws.Cells.Item(5,7).Value = "Application fees:https://www.ucl.ac.uk/prospective-students/graduate/application-fees"
$ws.Cells.Item(5,8).Value = "Entry requirements:#entry-requirements"
$ws.Cells.Item(5,9).Value = "Graduate degrees:/prospective-students/graduate/graduate-degrees"
$ws.Cells.Item(5,10).Value = "Taught Degrees:/prospective-students/graduate/taught-degrees"
$ws.Cells.Item(5,11).Value = "Applying for Graduate Taught Study at UCL:/prospective-students/graduate/taught-degrees/applying-graduate-taught-study-ucl"
$ws.Cells.Item(5,12).Value = "Research Degrees:/prospective-students/graduate/research-degrees"
$ws.Cells.Item(5,13).Value = "Applying for Graduate Research Study at UCL:/prospective-students/graduate/research-degrees/applying-graduate-research-study-ucl"
$ws.Cells.Item(5,14).Value = "Entry requirements:/prospective-students/graduate/teacher-training/entry-requirements"
$ws.Cells.Item(5,15).Value = "How to apply:/prospective-students/graduate/teacher-training/how-apply"

# Row 6
$ws.Cells.Item(6,1).Value = "advanced audiology: otology and skull base surgery msc"
$ws.Cells.Item(6,2).Value = "https://www.ucl.ac.uk/prospective-students/graduate/taught-degrees/advanced-audiology-otology-and-skull-base-surgery-msc"
$ws.Cells.Item(6,3).Value = "Faculty of Brain Sciences | Ear Institute"
$ws.Cells.Item(6,4).Value = "This MSc is structured to provide the theoretical, scientific, clinical, research and vocational skills necessary to practise enhanced otology, audiovestibular medicine and audiology. The programme is designed for ENT surgeons, audiovestibular physicians, paediatricians, GPs, neurologists and other trainees and physicians with an interest in the medical aspects of audiology who wish to develop or enhance their careers or specialise in otology."
$ws.Cells.Item(6,5).Value = "English language requirements:https://www.ucl.ac.uk/prospective-students/graduate/learning-and-living-ucl/international-students/english-language-requirements"
$ws.Cells.Item(6,6).Value = "UCL Graduate Admissions team:https://www.ucl.ac.uk/prospective-students/graduate/admissions-enquiries#form"
$ws.Cells.Item(6,7).Value = "Application fees:https://www.ucl.ac.uk/prospective-students/graduate/application-fees"
$ws.Cells.Item(6,8).Value = "Entry requirements:#entry-requirements"
$ws.Cells.Item(6,9).Value = "Graduate degrees:/prospective-students/graduate/graduate-degrees"
$ws.Cells.Item(6,10).Value = "Taught Degrees:/prospective-students/graduate/taught-degrees"
$ws.Cells.Item(6,11).Value = "Applying for Graduate Taught Study at UCL:/prospective-students/graduate/taught-degrees/applying-graduate-taught-study-ucl"
$ws.Cells.Item(6,12).Value = "Research Degrees:/prospective-students/graduate/research-degrees"
$ws.Cells.Item(6,13).Value = "Applying for Graduate Research Study at UCL:/prospective-students/graduate/research-degrees/applying-graduate-research-study-ucl"
$ws.Cells.Item(6,14).Value = "Entry requirements:/prospective-students/graduate/teacher-training/entry-requirements"
$ws.Cells.Item(6,15).Value = "How to apply:/prospective-students/graduate/teacher-training/how-apply"

# Row 7
$ws.Cells.Item(7,1).Value = "advanced biomedical imaging msc"
$ws.Cells.Item(7,2).Value = "https://www.ucl.ac.uk/prospective-students/graduate/taught-degrees/advanced-biomedical-imaging-msc"
$ws.Cells.Item(7,3).Value = "Faculty of Medical Sciences | Division of Medicine"
$ws.Cells.Item(7,4).Value = "Imaging has contributed to some of the most significant advances in biomedicine and healthcare. This one-year MSc will equip you with detailed knowledge of the imaging techniques that are shaping biomedical research and translational medicine around the world. You will develop the skills and research experience needed to progress in this highly sought-after field."
$ws.Cells.Item(7,5).Value = "English language requirements:https://www.ucl.ac.uk/prospective-students/graduate/learning-and-living-ucl/international-students/english-language-requirements"
$ws.Cells.Item(7,6).Value = "Application fees:https://www.ucl.ac.uk/prospective-students/graduate/application-fees"
$ws.Cells.Item(7,7).Value = "Entry requirements:#entry-requirements"
$ws.Cells.Item(7,8).Value = "Graduate degrees:/prospective-students/graduate/graduate-degrees"
$ws.Cells.Item(7,9).Value = "Taught Degrees:/prospective-students/graduate/taught-degrees"
$ws.Cells.Item(7,10).Value = "Applying for Graduate Taught Study at UCL:/prospective-students/graduate/taught-degrees/applying-graduate-taught-study-ucl"
$ws.Cells.Item(7,11).Value = "Research Degrees:/prospective-students/graduate/research-degrees"
$ws.Cells.Item(7,12).Value = "Applying for Graduate Research Study at UCL:/prospective-students/graduate/research-degrees/applying-graduate-research-study-ucl"
$ws.Cells.Item(7,13).Value = "Entry requirements:/prospective-students/graduate/teacher-training/entry-requirements"
$ws.Cells.Item(7,14).Value = "How to apply:/prospective-students/graduate/teacher-training/how-apply"

# Row 8
$ws.Cells.Item(8,1).Value = "advanced clinical practice in ophthalmology (integrated degree apprenticeship) msc"
$ws.Cells.Item(8,2).Value = "https://www.ucl.ac.uk/prospective-students/graduate/taught-degrees/advanced-clinical-practice-ophthalmology-integrated-degree-apprenticeship-msc"
$ws.Cells.Item(8,3).Value = "Faculty of Brain Sciences | Institute of Ophthalmology"
$ws.Cells.Item(8,4).Value = "The aim of this programme is to develop clinical practice and to contribute to the advancement of eye care delivery in primary, secondary or tertiary care settings. The programme has been developed by practicing clinicians for multidisciplinary ophthalmic non-medical clinicians such as optometrists, orthoptists and ophthalmic nurses. The structure of the programme is based on the four pillars of clinical, leadership, education and research."
$ws.Cells.Item(8,5).Value = "English language requirements:https://www.ucl.ac.uk/prospective-students/graduate/learning-and-living-ucl/international-students/english-language-requirements"
$ws.Cells.Item(8,6).Value = "ioo.admissions@ucl.ac.uk:mailto:ioo.admissions@ucl.ac.uk"
$ws.Cells.Item(8,7).Value = "ioo.admissions@ucl.ac.uk:mailto:ioo.admissions@ucl.ac.uk"
$ws.Cells.Item(8,8).Value = "Entry requirements:#entry-requirements"
$ws.Cells.Item(8,9).Value = "Graduate degrees:/prospective-students/graduate/graduate-degrees"
$ws.Cells.Item(8,10).Value = "Taught Degrees:/prospective-students/graduate/taught-degrees"
$ws.Cells.Item(8,11).Value = "Applying for Graduate Taught Study at UCL:/prospective-students/graduate/taught-degrees/applying-graduate-taught-study-ucl"
$ws.Cells.Item(8,12).Value = "Research Degrees:/prospective-students/graduate/research-degrees"
$ws.Cells.Item(8,13).Value = "Applying for Graduate Research Study at UCL:/prospective-students/graduate/research-degrees/applying-graduate-research-study-ucl"
$ws.Cells.Item(8,14).Value = "Entry requirements:/prospective-students/graduate/teacher-training/entry-requirements"
$ws.Cells.Item(8,15).Value = "How to apply:/prospective-students/graduate/teacher-training/how-apply"

# Row 9
$ws.Cells.Item(9,1).Value = "advanced educational practice grad dip"
$ws.Cells.Item(9,2).Value = "https://www.ucl.ac.uk/prospective-students/graduate/taught-degrees/advanced-educational-practice-grad-dip"
$ws.Cells.Item(9,3).Value = "IOE | Curriculum, Pedagogy and Assessment"
$ws.Cells.Item(9,4).Value = "This Graduate Diploma has been developed to enable participants to reflect on their own educational practice, taking a professional academic approach. It works as excellent continuing professional development for teachers, teaching assistants and others working in educational settings. The flexibility is ideal for those working full-time, such that modules can be studied in any order. "
$ws.Cells.Item(9,5).Value = "Apply for this course:#programme-choice"
$ws.Cells.Item(9,6).Value = "English language requirements:https://www.ucl.ac.uk/prospective-students/graduate/learning-and-living-ucl/international-students/english-language-requirements"
$ws.Cells.Item(9,7).Value = "Application fees:https://www.ucl.ac.uk/prospective-students/graduate/application-fees"
$ws.Cells.Item(9,8).Value = "Application Guidance:https://www.ucl.ac.uk/prospective-students/graduate/apply"
$ws.Cells.Item(9,9).Value = "Apply for this course:#"
$ws.Cells.Item(9,10).Value = "Entry requirements:#entry-requirements"
$ws.Cells.Item(9,11).Value = "Graduate degrees:/prospective-students/graduate/graduate-degrees"
$ws.Cells.Item(9,12).Value = "Taught Degrees:/prospective-students/graduate/taught-degrees"
$ws.Cells.Item(9,13).Value = "Applying for Graduate Taught Study at UCL:/prospective-students/graduate/taught-degrees/applying-graduate-taught-study-ucl"
$ws.Cells.Item(9,14).Value = "Research Degrees:/prospective-students/graduate/research-degrees"
$ws.Cells.Item(9,15).Value = "Applying for Graduate Research Study at UCL:/prospective-students/graduate/research-degrees/applying-graduate-research-study-ucl"
$ws.Cells.Item(9,16).Value = "Entry requirements:/prospective-students/graduate/teacher-training/entry-requirements"
$ws.Cells.Item(9,17).Value = "How to apply:/prospective-students/graduate/teacher-training/how-apply"

# Row 10
$ws.Cells.Item(10,1).Value = "advanced materials science (data-driven innovation) msc"
$ws.Cells.Item(10,2).Value = "https://www.ucl.ac.uk/prospective-students/graduate/taught-degrees/advanced-materials-science-data-driven-innovation-msc"
$ws.Cells.Item(10,3).Value = "Faculty of Mathematical and Physical Sciences | Faculty of Mathematical and Physical Sciences"
$ws.Cells.Item(10,4).Value = "The digital revolution and recent advents in data science, machine learning (ML) and artificial intelligence (AI) have sparked demand for next generation materials data scientists, able to utilise these emerging technologies for enhanced materials design and discovery. This programme will enable students to explore how the establishment of  Processing-Structure-Properties-Performance (PSPP) relationships can be significantly enhanced using data driven approaches."
$ws.Cells.Item(10,5).Value = "English language requirements:https://www.ucl.ac.uk/prospective-students/graduate/learning-and-living-ucl/international-students/english-language-requirements"
$ws.Cells.Item(10,6).Value = "Application fees:https://www.ucl.ac.uk/prospective-students/graduate/application-fees"
$ws.Cells.Item(10,7).Value = "Entry requirements:#entry-requirements"
$ws.Cells.Item(10,8).Value = "Graduate degrees:/prospective-students/graduate/graduate-degrees"
$ws.Cells.Item(10,9).Value = "Taught Degrees:/prospective-students/graduate/taught-degrees"
$ws.Cells.Item(10,10).Value = "Applying for Graduate Taught Study at UCL:/prospective-students/graduate/taught-degrees/applying-graduate-taught-study-ucl"
$ws.Cells.Item(10,11).Value = "Research Degrees:/prospective-students/graduate/research-degrees"
$ws.Cells.Item(10,12).Value = "Applying for Graduate Research Study at UCL:/prospective-students/graduate/research-degrees/applying-graduate-research-study-ucl"
$ws.Cells.Item(10,13).Value = "Entry requirements:/prospective-students/graduate/teacher-training/entry-requirements"
$ws.Cells.Item(10,14).Value = "How to apply:/prospective-students/graduate/teacher-training/how-apply"

# Row 11
$ws.Cells.Item(11,1).Value = "advanced materials science (energy storage) msc"
$ws.Cells.Item(11,2).Value = "https://www.ucl.ac.uk/prospective-students/graduate/taught-degrees/advanced-materials-science-energy-storage-msc"
$ws.Cells.Item(11,3).Value = "Faculty of Mathematical and Physical Sciences | Faculty of Mathematical and Physical Sciences"
$ws.Cells.Item(11,4).Value = "With global challenges in climate, environment, healthcare and economy demand, there is increasing need for scientific experts and entrepreneurs who can develop novel materials with advanced properties - addressing critical issues from energy to healthcare - and take scientific discoveries to the commercial world. This degree combines frontline research-based teaching from across UCL to train the next generation of materials scientists."
$ws.Cells.Item(11,5).Value = "English language requirements:https://www.ucl.ac.uk/prospective-students/graduate/learning-and-living-ucl/international-students/english-language-requirements"
$ws.Cells.Item(11,6).Value = "Application fees:https://www.ucl.ac.uk/prospective-students/graduate/application-fees"
$ws.Cells.Item(11,7).Value = "Entry requirements:#entry-requirements"
$ws.Cells.Item(11,8).Value = "Graduate degrees:/prospective-students/graduate/graduate-degrees"
$ws.Cells.Item(11,9).Value = "Taught Degrees:/prospective-students/graduate/taught-degrees"
$ws.Cells.Item(11,10).Value = "Applying for Graduate Taught Study at UCL:/prospective-students/graduate/taught-degrees/applying-graduate-taught-study-ucl"
$ws.Cells.Item(11,11).Value = "Research Degrees:/prospective-students/graduate/research-degrees"
$ws.Cells.Item(11,12).Value = "Applying for Graduate Research Study at UCL:/prospective-students/graduate/research-degrees/applying-graduate-research-study-ucl"
$ws.Cells.Item(11,13).Value = "Entry requirements:/prospective-students/graduate/teacher-training/entry-requirements"
$ws.Cells.Item(11,14).Value = "How to apply:/prospective-students/graduate/teacher-training/how-apply"

